# --- jobactions.xlsx update: add job-action entries from 2025-03-15 job alerts ---
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Mark "No" (no response) in column E for three earlier rows the author revisited
$ws.Cells.Item(94, 5).Value = "No"
$ws.Cells.Item(98, 5).Value = "No"
$ws.Cells.Item(100, 5).Value = "No"

# New job-action rows 112-118 (date serial 45731 = 2025-03-15)
$newRowDate = 45731

$ws.Cells.Item(112, 1).Value = $newRowDate
$ws.Cells.Item(111, 1).Copy()
$ws.Cells.Item(112, 1).PasteSpecial(-4122)
$ws.Cells.Item(112, 2).Value = 'Amazon'
$ws.Cells.Item(112, 3).Value = 'Principal Data Scientist, Forecasting, ASIN Forecasting'
$ws.Cells.Item(112, 4).Value = 'demand forecasting'
$ws.Cells.Item(112, 6).Value = 'https://www.linkedin.com/jobs/view/4137885098/?trackingId=e299vDboTfLyewOezp%2F49Q%3D%3D&refId=ByteString%28length%3D16%2Cbytes%3D39f6ef31...df0aa05c%29&midToken=AQFpnZsm4rTQjw&midSig=0b_0l8pijg0XI1&trkEmail=eml-email_job_alert_digest_01-job_card-0-jobcard_body-null-1j75g~m89l2nn0~g1-null-null&eid=1j75g-m89l2nn0-g1&otpToken=MTAwMDE5ZTUxMTJhYzFjZWI1MjkwMWU4NDYxZWUyYjA4OWNkZDU0NDliYTQ4NzZmNzdjMTAwNmU0NzViNTY5NThlOWZiN2EzNjhmYWYzOTFhN2JiZmU2OTkzNjc3ZjdlYjYzNWJjY2FhYTY1YTAsMSwx'

$ws.Cells.Item(113, 1).Value = $newRowDate
$ws.Cells.Item(111, 1).Copy()
$ws.Cells.Item(113, 1).PasteSpecial(-4122)
$ws.Cells.Item(113, 2).Value = 'Amazon'
$ws.Cells.Item(113, 3).Value = 'Sr. Data Scientist, Perceptor (Kumo Analytics), AWS Support'
$ws.Cells.Item(113, 4).Value = 'business stuff?  They mention forecasting'
$ws.Cells.Item(113, 6).Value = 'https://www.linkedin.com/jobs/view/4184889780/?trackingId=8RyyZELEcPLpz3qE5hd%2BFA%3D%3D&refId=ByteString%28length%3D16%2Cbytes%3D39f6ef31...df0aa05c%29&midToken=AQFpnZsm4rTQjw&midSig=0b_0l8pijg0XI1&trkEmail=eml-email_job_alert_digest_01-job_card-0-jobcard_body-null-1j75g~m89l2nn0~g1-null-null&eid=1j75g-m89l2nn0-g1&otpToken=MTAwMDE5ZTUxMTJhYzFjZWI1MjkwMWU4NDYxZWUyYjA4OWNkZDU0NDliYTQ4NzZmNzdjMTAwNmU0NzViNTY5NThlOWZiN2EzNjhmYWYzOTFhN2JiZmU2OTkzNjc3ZjdlYjYzNWJjY2FhYTY1YTAsMSwx'

$ws.Cells.Item(114, 1).Value = $newRowDate
$ws.Cells.Item(111, 1).Copy()
$ws.Cells.Item(114, 1).PasteSpecial(-4122)
$ws.Cells.Item(114, 2).Value = 'Amazon'
$ws.Cells.Item(114, 3).Value = 'Senior Data Scientist, Last Mile Science'
$ws.Cells.Item(114, 4).Value = 'logistics?  Mention forecasting'
$ws.Cells.Item(114, 6).Value = 'https://www.linkedin.com/jobs/view/4184242488/?trackingId=EQII1MqLI2PsUdCmISGePw%3D%3D&refId=ByteString%28length%3D16%2Cbytes%3D002c1553...c1ec6cde%29&midToken=AQFpnZsm4rTQjw&midSig=0YR_8CJIUw_HE1&trkEmail=eml-email_job_alert_digest_01-job_card-0-jobcard_body-null-1j75g~m88528ge~mh-null-null&eid=1j75g-m88528ge-mh&otpToken=MTAwMDE5ZTUxMTJhYzFjZWI1MjkwMWViNGYxZmVmYjU4ZmNiZDM0OTllYTQ4NzZmNzdjMTAwNmU0NzViNTY4MWFhYWViNDgzMTNiOGIxMzY1NWE1NTY3ZDA1NzdkNmM0NjcyMTI5OWIyYmM3NzcsMSwx'

$ws.Cells.Item(115, 1).Value = $newRowDate
$ws.Cells.Item(111, 1).Copy()
$ws.Cells.Item(115, 1).PasteSpecial(-4122)
$ws.Cells.Item(115, 2).Value = 'Amazon'
$ws.Cells.Item(115, 3).Value = 'Data Scientist, Topline Forecasting'
$ws.Cells.Item(115, 4).Value = 'Biz forecasting'
$ws.Cells.Item(115, 6).Value = 'https://www.linkedin.com/jobs/view/4007093976/?trackingId=GXAgJQcAHrymvJVwxCcwfQ%3D%3D&refId=ByteString%28length%3D16%2Cbytes%3D96b02801...3f23e2c7%29&midToken=AQFpnZsm4rTQjw&midSig=3GhlV4E8MsYXE1&trkEmail=eml-email_job_alert_digest_01-job_card-0-jobcard_body-null-1j75g~m85cbosa~6t-null-null&eid=1j75g-m85cbosa-6t&otpToken=MTAwMDE5ZTUxMTJhYzFjZWI1MjkwMWViNDExYWVmYmQ4ZWNkZDY0NzlmYTQ4NzZmNzdjMTAwNmU0NzViNTY5OTgyZDRhMTk5NTRmMWZmZTFlMjE1NzI2M2E1MWU0ODBmYjQ3OTI3ZTg2ODA5NjMsMSwx'

$ws.Cells.Item(116, 1).Value = $newRowDate
$ws.Cells.Item(111, 1).Copy()
$ws.Cells.Item(116, 1).PasteSpecial(-4122)
$ws.Cells.Item(116, 2).Value = 'Amazon'
$ws.Cells.Item(116, 3).Value = 'Sr. Data Scientist, Devices Decision Scienc'
$ws.Cells.Item(116, 4).Value = 'biz stuff, mention forecasting'
$ws.Cells.Item(116, 6).Value = 'https://www.linkedin.com/jobs/view/4148960536/?trackingId=zHv7uEAnsa0%2FK6LDHqPWow%3D%3D&refId=ByteString%28length%3D16%2Cbytes%3D4987540d...389951e1%29&midToken=AQFpnZsm4rTQjw&midSig=2wsph5SHHgRHE1&trkEmail=eml-email_job_alert_digest_01-job_card-0-jobcard_body-null-1j75g~m7y5jufr~ei-null-null&eid=1j75g-m7y5jufr-ei&otpToken=MTAwMDE5ZTUxMTJhYzFjZWI1MjkwMWViNDUxZmUzYjY4OWNhZDU0NjkwYTQ4NzZmNzdjMTAwNmU0NzViNTZiOWY0OTNhMTkwNTRlNWZiNDVkOTEyYTY5ZGYwNWJkOWEwNzRiYjg5MGQ2OTAxNTIsMSwx'

$ws.Cells.Item(117, 1).Value = $newRowDate
$ws.Cells.Item(111, 1).Copy()
$ws.Cells.Item(117, 1).PasteSpecial(-4122)
$ws.Cells.Item(117, 2).Value = 'Amazon'
$ws.Cells.Item(117, 3).Value = 'Sr. Applied Scientist, Renewable Energy Optimization'
$ws.Cells.Item(117, 4).Value = 'RES opt but want forecasting experience'
$ws.Cells.Item(117, 6).Value = 'https://www.amazon.jobs/en/jobs/2913322/sr-applied-scientist-renewable-energy-optimization?cmpid=DA_INAD200785B'

$ws.Cells.Item(118, 1).Value = $newRowDate
$ws.Cells.Item(111, 1).Copy()
$ws.Cells.Item(118, 1).PasteSpecial(-4122)
$ws.Cells.Item(118, 2).Value = 'Strella'
$ws.Cells.Item(118, 3).Value = 'Data Scientist'
$ws.Cells.Item(118, 4).Value = 'produce transport decisions, I think, Seattle'
$ws.Cells.Item(118, 6).Value = 'https://www.linkedin.com/jobs/view/4175135972/?refId=ByteString(length%3D16%2Cbytes%3D8a8662a5...2cc15d02)&trackingId=Ib7%2BrwAH6ArYZYJILpsqvg%3D%3D'

$excel.CutCopyMode = 0

# Leave the selection where the author last edited (E100)
$ws.Range("E100").Select()

Write-Host "done"
